# chore: update Sheets via scheduled runner
# Refreshes the derived price/profit columns (currentAveragePrice,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ) on the
# per-job leve-profit tables with newly scraped market-board figures.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 114.5
$ws.Range("I2").Value = 114.5
$ws.Range("K2").Value = 114.5
$ws.Range("M2").Value = -1.5
$ws.Range("H11").Value = 69.09524
$ws.Range("I11").Value = 69.09524
$ws.Range("K11").Value = 69.09524
$ws.Range("M11").Value = 70.90476
$ws.Range("H17").Value = 2508.923
$ws.Range("J17").Value = 1462.75
$ws.Range("L17").Value = 4388.25
$ws.Range("N17").Value = -4724.25
$ws.Range("H19").Value = 826.2308
$ws.Range("J19").Value = 749
$ws.Range("L19").Value = 749
$ws.Range("N19").Value = -1099
$ws.Range("H62").Value = 3999.5
$ws.Range("I62").Value = 3999.5
$ws.Range("K62").Value = 3999.5
$ws.Range("M62").Value = -3375.5
$ws.Range("H65").Value = 3999.5
$ws.Range("I65").Value = 3999.5
$ws.Range("K65").Value = 19997.5
$ws.Range("M65").Value = -16877.5
$ws.Range("H116").Value = 4088.9092
$ws.Range("I116").Value = 4127.9
$ws.Range("K116").Value = 4127.9
$ws.Range("M116").Value = -685.8999999999996
$ws.Range("H131").Value = 14300
$ws.Range("I131").Value = 15750
$ws.Range("J131").Value = 13333.333
$ws.Range("K131").Value = 47250
$ws.Range("L131").Value = 39999.999
$ws.Range("M131").Value = -42210
$ws.Range("N131").Value = -50079.999
$ws.Range("H138").Value = 3251.0571
$ws.Range("J138").Value = 4531.727
$ws.Range("L138").Value = 13595.181
$ws.Range("N138").Value = -23875.181

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10754947
$ws.Range("J32").Value = 2266
$ws.Range("L32").Value = 2266
$ws.Range("N32").Value = -2840
$ws.Range("H45").Value = 2746.889
$ws.Range("I45").Value = 2982.4
$ws.Range("K45").Value = 2982.4
$ws.Range("M45").Value = -2605.4
$ws.Range("H74").Value = 2788.9285
$ws.Range("I74").Value = 2398.4211
$ws.Range("J74").Value = 6498.75
$ws.Range("K74").Value = 2398.4211
$ws.Range("L74").Value = 6498.75
$ws.Range("M74").Value = -1524.4211
$ws.Range("N74").Value = -8246.75
$ws.Range("H77").Value = 2788.9285
$ws.Range("I77").Value = 2398.4211
$ws.Range("J77").Value = 6498.75
$ws.Range("K77").Value = 11992.1055
$ws.Range("L77").Value = 32493.75
$ws.Range("M77").Value = -7624.1055
$ws.Range("N77").Value = -41229.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2024.1666
$ws.Range("I132").Value = 2006.1904
$ws.Range("J132").Value = 2150
$ws.Range("K132").Value = 6018.5712
$ws.Range("L132").Value = 6450
$ws.Range("M132").Value = -3488.5712
$ws.Range("N132").Value = -11510
$ws.Range("H134").Value = 2020.3778
$ws.Range("I134").Value = 1995.6428
$ws.Range("K134").Value = 5986.928400000001
$ws.Range("M134").Value = -3451.928400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1448.9474
$ws.Range("I5").Value = 1444.5
$ws.Range("J5").Value = 1472.6666
$ws.Range("K5").Value = 4333.5
$ws.Range("L5").Value = 4417.9998
$ws.Range("M5").Value = -4221.5
$ws.Range("N5").Value = -4641.9998
$ws.Range("H25").Value = 175
$ws.Range("I25").Value = 177
$ws.Range("J25").Value = 170
$ws.Range("K25").Value = 531
$ws.Range("L25").Value = 510
$ws.Range("M25").Value = -362
$ws.Range("N25").Value = -848
$ws.Range("H30").Value = 175
$ws.Range("I30").Value = 177
$ws.Range("J30").Value = 170
$ws.Range("K30").Value = 531
$ws.Range("L30").Value = 510
$ws.Range("M30").Value = -429
$ws.Range("N30").Value = -714
$ws.Range("H48").Value = 792.3077
$ws.Range("I48").Value = 550
$ws.Range("J48").Value = 900
$ws.Range("K48").Value = 1650
$ws.Range("L48").Value = 2700
$ws.Range("M48").Value = -1400
$ws.Range("N48").Value = -3200
$ws.Range("H50").Value = 1216.6
$ws.Range("J50").Value = 1483.25
$ws.Range("L50").Value = 4449.75
$ws.Range("N50").Value = -5411.75
$ws.Range("H53").Value = 1216.6
$ws.Range("J53").Value = 1483.25
$ws.Range("L53").Value = 4449.75
$ws.Range("N53").Value = -5411.75
$ws.Range("H55").Value = 723035.7
$ws.Range("I55").Value = 833.3333
$ws.Range("J55").Value = 920000
$ws.Range("K55").Value = 2499.9999
$ws.Range("L55").Value = 2760000
$ws.Range("M55").Value = -2322.9999
$ws.Range("N55").Value = -2760354
$ws.Range("H112").Value = 0
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("M112").ClearContents()
$ws.Range("N112").ClearContents()
$ws.Range("H113").Value = 1641.8334
$ws.Range("J113").Value = 1497.8
$ws.Range("L113").Value = 4493.4
$ws.Range("N113").Value = -8833.4
$ws.Range("H131").Value = 3642.2144
$ws.Range("J131").Value = 6329.4
$ws.Range("L131").Value = 18988.2
$ws.Range("N131").Value = -29068.2
$ws.Range("H135").Value = 1448.9474
$ws.Range("I135").Value = 1444.5
$ws.Range("J135").Value = 1472.6666
$ws.Range("K135").Value = 13000.5
$ws.Range("L135").Value = 13253.9994
$ws.Range("M135").Value = -10465.5
$ws.Range("N135").Value = -18323.9994
$ws.Range("H137").Value = 2817.5386
$ws.Range("I137").Value = 2466.5
$ws.Range("J137").Value = 3118.4285
$ws.Range("K137").Value = 7399.5
$ws.Range("L137").Value = 9355.2855
$ws.Range("M137").Value = -2299.5
$ws.Range("N137").Value = -19555.2855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4362.231
$ws.Range("I46").Value = 1900
$ws.Range("J46").Value = 5100.9
$ws.Range("K46").Value = 1900
$ws.Range("L46").Value = 5100.9
$ws.Range("M46").Value = -1712
$ws.Range("N46").Value = -5476.9
$ws.Range("H93").Value = 3190.5334
$ws.Range("I93").Value = 3987.625
$ws.Range("J93").Value = 2279.5715
$ws.Range("K93").Value = 3987.625
$ws.Range("L93").Value = 2279.5715
$ws.Range("M93").Value = -2739.625
$ws.Range("N93").Value = -4775.5715

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 736.8421
$ws.Range("I113").Value = 462.6154
$ws.Range("J113").Value = 1331
$ws.Range("K113").Value = 1387.8462
$ws.Range("L113").Value = 3993
$ws.Range("M113").Value = 782.1538
$ws.Range("N113").Value = -8333
